# Implements "Added implementation of MSM measure." for the
# mall-search_structure workbook: the "interfaceOperations" sheet gains
# rows for the inherited java.lang.Object operations (equals, toString,
# getClass, notifyAll, hashCode, wait, notify, wait(long), wait(long,int))
# interleaved with the pre-existing EsProductController operations, for a
# total of 18 data rows (was 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interfaceOperations")

$interfaceName = "com.macro.mall.search.controller.EsProductController"
$modifier = "public"

# Final row order (Operation Signature, Return Type) for rows 2..19
$rows = @(
    @("equals(java.lang.Object)", "boolean"),
    @("toString()", "java.lang.String"),
    @("importAllList()", "com.macro.mall.common.api.CommonResult"),
    @("search(java.lang.String, java.lang.Integer, java.lang.Integer)", "com.macro.mall.common.api.CommonResult"),
    @("search(java.lang.String, java.lang.Long, java.lang.Long, java.lang.Integer, java.lang.Integer, java.lang.Integer)", "com.macro.mall.common.api.CommonResult"),
    @("recommend(java.lang.Long, java.lang.Integer, java.lang.Integer)", "com.macro.mall.common.api.CommonResult"),
    @("getClass()", "java.lang.Class"),
    @("notifyAll()", "void"),
    @("hashCode()", "int"),
    @("wait()", "void"),
    @("searchRelatedInfo(java.lang.String)", "com.macro.mall.common.api.CommonResult"),
    @("EsProductController()", "void"),
    @("create(java.lang.Long)", "com.macro.mall.common.api.CommonResult"),
    @("notify()", "void"),
    @("delete(java.util.List)", "com.macro.mall.common.api.CommonResult"),
    @("wait(long)", "void"),
    @("wait(long, int)", "void"),
    @("delete(java.lang.Long)", "com.macro.mall.common.api.CommonResult")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $interfaceName
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $modifier
    $ws.Cells.Item($r, 4).Value = $row[1]
    $r = $r + 1
}
